{"js": "// Replace each two-digit-by-two-digit multiplication prompt in the worksheet\n// table with a newly generated equation, in document order. The table has\n// many empty \"answer\" rows interleaved with rows that hold the \"NN\u00d7NN=\" text;\n// only the non-empty cells are targeted, exactly as in the source diff.\nconst replacements = [\n  { oldVal: \"79\u00d767=\", newVal: \"64\u00d781=\" },\n  { oldVal: \"52\u00d766=\", newVal: \"19\u00d714=\" },\n  { oldVal: \"70\u00d744=\", newVal: \"59\u00d780=\" },\n  { oldVal: \"34\u00d750=\", newVal: \"65\u00d725=\" },\n  { oldVal: \"70\u00d736=\", newVal: \"17\u00d796=\" },\n  { oldVal: \"14\u00d725=\", newVal: \"15\u00d736=\" },\n  { oldVal: \"33\u00d779=\", newVal: \"61\u00d749=\" },\n  { oldVal: \"20\u00d747=\", newVal: \"52\u00d782=\" },\n  { oldVal: \"32\u00d757=\", newVal: \"42\u00d730=\" },\n  { oldVal: \"12\u00d790=\", newVal: \"20\u00d711=\" },\n  { oldVal: \"19\u00d724=\", newVal: \"62\u00d746=\" },\n  { oldVal: \"76\u00d724=\", newVal: \"13\u00d759=\" },\n  { oldVal: \"36\u00d780=\", newVal: \"84\u00d792=\" },\n  { oldVal: \"97\u00d788=\", newVal: \"98\u00d745=\" },\n  { oldVal: \"65\u00d725=\", newVal: \"57\u00d739=\" },\n  { oldVal: \"35\u00d791=\", newVal: \"59\u00d757=\" },\n  { oldVal: \"72\u00d753=\", newVal: \"80\u00d765=\" },\n  { oldVal: \"81\u00d740=\", newVal: \"20\u00d793=\" },\n  { oldVal: \"94\u00d796=\", newVal: \"22\u00d761=\" },\n  { oldVal: \"63\u00d721=\", newVal: \"67\u00d726=\" },\n  { oldVal: \"45\u00d733=\", newVal: \"14\u00d721=\" },\n  { oldVal: \"23\u00d777=\", newVal: \"80\u00d797=\" },\n  { oldVal: \"71\u00d790=\", newVal: \"91\u00d751=\" },\n  { oldVal: \"77\u00d722=\", newVal: \"94\u00d784=\" },\n  { oldVal: \"53\u00d799=\", newVal: \"28\u00d746=\" }\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst grid = table.values;\nlet next = 0;\nfor (let r = 0; r < grid.length; r++) {\n  for (let c = 0; c < grid[r].length; c++) {\n    const text = grid[r][c];\n    if (text === \"\" || text === undefined || text === null) continue;\n    if (next >= replacements.length) continue;\n    const cell = table.getCell(r, c);\n    cell.value = replacements[next].newVal;\n    next++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each two-digit-by-two-digit multiplication prompt in the worksheet\n# table with a newly generated equation, in document order. The table has\n# many empty \"answer\" rows interleaved with rows that hold the \"NN\u00d7NN=\" text;\n# only the non-empty cells are targeted, exactly as in the source diff.\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n$newValues = @(\n  \"64\u00d781=\",\n  \"19\u00d714=\",\n  \"59\u00d780=\",\n  \"65\u00d725=\",\n  \"17\u00d796=\",\n  \"15\u00d736=\",\n  \"61\u00d749=\",\n  \"52\u00d782=\",\n  \"42\u00d730=\",\n  \"20\u00d711=\",\n  \"62\u00d746=\",\n  \"13\u00d759=\",\n  \"84\u00d792=\",\n  \"98\u00d745=\",\n  \"57\u00d739=\",\n  \"59\u00d757=\",\n  \"80\u00d765=\",\n  \"20\u00d793=\",\n  \"22\u00d761=\",\n  \"67\u00d726=\",\n  \"14\u00d721=\",\n  \"80\u00d797=\",\n  \"91\u00d751=\",\n  \"94\u00d784=\",\n  \"28\u00d746=\"\n)\n\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n$next = 0\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n  for ($c = 1; $c -le $colCount; $c++) {\n    $cell = $t.Cell($r, $c)\n    $text = $cell.Range.Text\n    if ($text.Length -gt 2) {\n      if ($next -lt $newValues.Count) {\n        $cell.Range.Text = $newValues[$next]\n        $next = $next + 1\n      }\n    }\n  }\n}\n"}
